$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5382891893386841
$ws.Range("B1").Value = 0.7993694543838501
$ws.Range("C1").Value = 5.069862365722656
$ws.Range("D1").Value = 2.19356894493103
$ws.Range("E1").Value = 1.184125065803528
